$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new data row for 2022-Q3 at the top
#    of the data block (row 2), pushing the existing quarters down by one
#    row. The "A" index column (0,1,2,...) is positional and stays as-is
#    except for the brand-new last row, which needs a new "7".
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B8:D8").Copy($summary.Range("B9:D9"))
$summary.Range("B7:D7").Copy($summary.Range("B8:D8"))
$summary.Range("B6:D6").Copy($summary.Range("B7:D7"))
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))

$summary.Range("A8").Copy($summary.Range("A9"))
$summary.Range("A9").Value = 7

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 2.88

# -----------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right after "总计", duplicating
#    the layout/formatting of the "2022-Q2" sheet (same A1:H8 shape) and
#    replacing its values with the 2022-Q3 holdings. The original
#    "2022-Q2" sheet (and every sheet after it) is left untouched content-
#    wise; it simply shifts one tab position to the right.
# -----------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy([System.Reflection.Missing]::Value, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q3.Range("B2").Value = "001838"
$q3.Range("C2").Value = "国投瑞银国家安全灵活配置混合"
$q3.Range("D2").Value = "27.40"
$q3.Range("E2").Value = "94.42"
$q3.Range("F2").Value = "7.84"
$q3.Range("G2").Value = "2.1482"
$q3.Range("H2").Value = 7

$q3.Range("B3").Value = "004139"
$q3.Range("C3").Value = "中邮军民融合灵活配置混合"
$q3.Range("D3").Value = "11.73"
$q3.Range("E3").Value = "88.46"
$q3.Range("F3").Value = "4.74"
$q3.Range("G3").Value = "0.5560"
$q3.Range("H3").Value = 10

$q3.Range("B4").Value = "011001"
$q3.Range("C4").Value = "中邮兴荣价值一年持有期混合"
$q3.Range("D4").Value = "5.15"
$q3.Range("E4").Value = "40.76"
$q3.Range("F4").Value = "2.31"
$q3.Range("G4").Value = "0.1190"
$q3.Range("H4").Value = 8

$q3.Range("B5").Value = "014781"
$q3.Range("C5").Value = "建信兴衡优选一年持有期混合A"
$q3.Range("D5").Value = "1.77"
$q3.Range("E5").Value = "46.74"
$q3.Range("F5").Value = "2.27"
$q3.Range("G5").Value = "0.0402"
$q3.Range("H5").Value = 8

$q3.Range("B6").Value = "014782"
$q3.Range("C6").Value = "建信兴衡优选一年持有期混合C"
$q3.Range("D6").Value = "0.70"
$q3.Range("E6").Value = "46.74"
$q3.Range("F6").Value = "2.27"
$q3.Range("G6").Value = "0.0159"
$q3.Range("H6").Value = 8

$q3.Range("B7").Value = "004840"
$q3.Range("C7").Value = "东兴品牌精选灵活配置混合A"
$q3.Range("D7").Value = "0.01"
$q3.Range("E7").Value = "89.84"
$q3.Range("F7").Value = "5.28"
$q3.Range("G7").Value = "0.0005"
$q3.Range("H7").Value = 7

$q3.Range("B8").Value = "006442"
$q3.Range("C8").Value = "东兴品牌精选灵活配置混合C"
$q3.Range("D8").Value = "0.00"
$q3.Range("E8").Value = "89.84"
$q3.Range("F8").Value = "5.28"
$q3.Range("G8").Value = 0
$q3.Range("H8").Value = 7
